$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header title (row 2, merged D2:G2) keeps its text; shared-string index
# shifts automatically once the other strings are rewritten below. ---
$ws.Range("D2").Value = "REPORTE POR USUARIO"

# --- Row 3: the F3 placeholder cell disappears entirely (no longer used). ---
$ws.Range("F3").Clear()

# --- Row 6 header labels change to reflect the new columns. ---
$ws.Range("C6").Value = "Tomo"
$ws.Range("D6").Value = "Movimiento"
$ws.Range("E6").Value = "Fecha Solicitud"
$ws.Range("F6").Value = "Fecha Devoluciòn"
$ws.Range("G6").Value = "Fecha Entrega"
$ws.Range("H6").Value = "Disponibilidad"

# Row 6 no longer forces a fixed (wrapped, 30pt) height - let it auto-size.
$ws.Rows.Item(6).AutoFit()

# --- Column widths adjust slightly now that "Fecha Entrega" was inserted. ---
$ws.Columns.Item(4).ColumnWidth = 24.8333333
$ws.Columns.Item(5).ColumnWidth = 20.5
$ws.Columns.Item(6).ColumnWidth = 19.5
$ws.Columns.Item(7).ColumnWidth = 19.5
$ws.Columns.Item(8).ColumnWidth = 20.6666667
